$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Read the current layout -------------------------------------------------
$ur = $ws.UsedRange
$lastRow = $ur.Rows.Count
$lastCol = $ur.Columns.Count

# Column B (index 2) is the first metric column; column A holds the row label.
$firstDataCol = 2

# Read existing headers (row 1), columns B..lastCol
$headers = @()
for ($c = $firstDataCol; $c -le $lastCol; $c++) {
    $headers += $ws.Cells.Item(1, $c).Value2
}

# Read existing data rows (label in col A, values keyed by header name)
$labels = @()
$rowData = @()
for ($r = 2; $r -le $lastRow; $r++) {
    $labels += $ws.Cells.Item($r, 1).Value2
    $vals = @{}
    for ($c = $firstDataCol; $c -le $lastCol; $c++) {
        $name = $headers[$c - $firstDataCol]
        $vals[$name] = $ws.Cells.Item($r, $c).Value2
    }
    $rowData += $vals
}

# --- Build the new column order ----------------------------------------------
# Every "base" metric (e.g. R2) is immediately followed by its "X" companion
# (e.g. R2X) if one exists among the current headers. R7 is missing its X
# companion (R7X) entirely, so it is inserted as a brand new column (value 0).
$baseNames = @()
$xNames = @{}
foreach ($h in $headers) {
    if ($h.EndsWith("X")) {
        $xNames[$h] = $true
    } else {
        $baseNames += $h
    }
}

$newHeaders = @()
foreach ($b in $baseNames) {
    $newHeaders += $b
    $xName = $b + "X"
    if ($xNames.ContainsKey($xName) -or $b -eq "R7") {
        $newHeaders += $xName
    }
}

# --- Write everything back out in the new order -------------------------------
$newLastCol = $firstDataCol + $newHeaders.Length - 1

# If the header row grew (new columns were inserted), give the new header
# cell(s) the same formatting as the existing header cells (bold, bordered,
# centered) by copying the format from the first header cell.
if ($newLastCol -gt $lastCol) {
    for ($c = $lastCol + 1; $c -le $newLastCol; $c++) {
        $ws.Cells.Item(1, $firstDataCol).Copy() | Out-Null
        $ws.Cells.Item(1, $c).PasteSpecial(-4122) | Out-Null
    }
}

for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $ws.Cells.Item(1, $firstDataCol + $i).Value = $newHeaders[$i]
}

for ($r = 2; $r -le $lastRow; $r++) {
    $vals = $rowData[$r - 2]
    for ($i = 0; $i -lt $newHeaders.Length; $i++) {
        $name = $newHeaders[$i]
        if ($vals.ContainsKey($name)) {
            $v = $vals[$name]
        } else {
            $v = 0
        }
        $ws.Cells.Item($r, $firstDataCol + $i).Value = $v
    }
}
